$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 75.94118
$ws.Range("J9").Value = 50.142857
$ws.Range("L9").Value = 50.142857
$ws.Range("N9").Value = -388.142857

$ws.Range("H17").Value = 4361.4907
$ws.Range("J17").Value = 4438.5557
$ws.Range("L17").Value = 13315.6671
$ws.Range("N17").Value = -13651.6671

$ws.Range("H33").Value = 5025.25
$ws.Range("I33").Value = 6104.615
$ws.Range("K33").Value = 6104.615
$ws.Range("M33").Value = -5875.615

$ws.Range("H41").Value = 14493470
$ws.Range("I41").Value = 561.7143
$ws.Range("J41").Value = 20834118
$ws.Range("K41").Value = 561.7143
$ws.Range("L41").Value = 20834118
$ws.Range("M41").Value = -121.7143
$ws.Range("N41").Value = -20834998

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H132").Value = 10528957
$ws.Range("I132").Value = 11238518
$ws.Range("K132").Value = 33715554
$ws.Range("M132").Value = -33713024

$ws.Range("H136").Value = 196889.75
$ws.Range("J136").Value = 196889.75
$ws.Range("L136").Value = 196889.75
$ws.Range("N136").Value = -207089.75

$ws.Range("H137").Value = 112350.94
$ws.Range("I137").Value = 198546.44
$ws.Range("J137").Value = 1528.1428
$ws.Range("K137").Value = 595639.3200000001
$ws.Range("L137").Value = 4584.428400000001
$ws.Range("M137").Value = -593089.3200000001
$ws.Range("N137").Value = -9684.428400000001

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H140").Value = 120000
$ws.Range("J140").Value = 120000
$ws.Range("L140").Value = 120000
$ws.Range("N140").Value = -130360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2377.1226
$ws.Range("I32").Value = 1460.125
$ws.Range("K32").Value = 1460.125
$ws.Range("M32").Value = -1173.125

$ws.Range("H45").Value = 5757036
$ws.Range("I45").Value = 7572526
$ws.Range("K45").Value = 7572526
$ws.Range("M45").Value = -7572149

$ws.Range("H102").Value = 2317688.5
$ws.Range("I102").Value = 2690793.5
$ws.Range("J102").Value = 4438
$ws.Range("K102").Value = 2690793.5
$ws.Range("L102").Value = 4438
$ws.Range("M102").Value = -2689171.5
$ws.Range("N102").Value = -7682

$ws.Range("H122").Value = 374361.03
$ws.Range("I122").Value = 2045.4359
$ws.Range("K122").Value = 6136.307699999999
$ws.Range("M122").Value = -3686.307699999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 7939933.5
$ws.Range("I99").Value = 11907827
$ws.Range("J99").Value = 4147.3335
$ws.Range("K99").Value = 11907827
$ws.Range("L99").Value = 4147.3335
$ws.Range("M99").Value = -11906329
$ws.Range("N99").Value = -7143.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 232.36363
$ws.Range("I7").Value = 76.57143000000001
$ws.Range("J7").Value = 505
$ws.Range("K7").Value = 76.57143000000001
$ws.Range("L7").Value = 505
$ws.Range("M7").Value = 36.42856999999999
$ws.Range("N7").Value = -731

$ws.Range("H31").Value = 2991.15
$ws.Range("I31").Value = 869.9722
$ws.Range("J31").Value = 4184.3125
$ws.Range("K31").Value = 869.9722
$ws.Range("L31").Value = 4184.3125
$ws.Range("M31").Value = -574.9722
$ws.Range("N31").Value = -4774.3125

$ws.Range("H34").Value = 2991.15
$ws.Range("I34").Value = 869.9722
$ws.Range("J34").Value = 4184.3125
$ws.Range("K34").Value = 869.9722
$ws.Range("L34").Value = 4184.3125
$ws.Range("M34").Value = -667.9722
$ws.Range("N34").Value = -4588.3125

$ws.Range("H58").Value = 3158.5454
$ws.Range("I58").Value = 2927.6667
$ws.Range("K58").Value = 2927.6667
$ws.Range("M58").Value = -2724.6667

$ws.Range("H99").Value = 3535.4
$ws.Range("J99").Value = 4027.6
$ws.Range("L99").Value = 4027.6
$ws.Range("N99").Value = -7023.6

$ws.Range("H109").Value = 57986.5
$ws.Range("J109").Value = 57986.5
$ws.Range("L109").Value = 57986.5
$ws.Range("N109").Value = -60066.5

$ws.Range("H122").Value = 4164.6665
$ws.Range("I122").Value = 3873.25
$ws.Range("J122").Value = 4747.5
$ws.Range("K122").Value = 11619.75
$ws.Range("L122").Value = 14242.5
$ws.Range("M122").Value = -9169.75
$ws.Range("N122").Value = -19142.5

$ws.Range("H126").Value = 3535.4
$ws.Range("J126").Value = 4027.6
$ws.Range("L126").Value = 12082.8
$ws.Range("N126").Value = -17022.8

$ws.Range("H132").Value = 34408.71
$ws.Range("I132").Value = 39009.48
$ws.Range("K132").Value = 117028.44
$ws.Range("M132").Value = -114498.44

$ws.Range("H134").Value = 7639.385
$ws.Range("I134").Value = 5310.6787
$ws.Range("K134").Value = 15932.0361
$ws.Range("M134").Value = -13397.0361

$ws.Range("H136").Value = 3158.5454
$ws.Range("I136").Value = 2927.6667
$ws.Range("K136").Value = 8783.000100000001
$ws.Range("M136").Value = -6233.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10909652
$ws.Range("I4").Value = 17358724
$ws.Range("J4").Value = 429912.38
$ws.Range("K4").Value = 52076172
$ws.Range("L4").Value = 1289737.14
$ws.Range("M4").Value = -52076060
$ws.Range("N4").Value = -1289961.14

$ws.Range("H56").Value = 8149.25
$ws.Range("I56").Value = 8149.25
$ws.Range("K56").Value = 8149.25
$ws.Range("M56").Value = -7619.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13338799
$ws.Range("I70").Value = 14291070
$ws.Range("K70").Value = 14291070
$ws.Range("M70").Value = -14290800

$ws.Range("H73").Value = 13338799
$ws.Range("I73").Value = 14291070
$ws.Range("K73").Value = 14291070
$ws.Range("M73").Value = -14290134

$ws.Range("H102").Value = 3317507.2
$ws.Range("I102").Value = 3970088.5
$ws.Range("K102").Value = 3970088.5
$ws.Range("M102").Value = -3968466.5

$ws.Range("H107").Value = 1042.4286
$ws.Range("I107").Value = 2164
$ws.Range("J107").Value = 201.25
$ws.Range("K107").Value = 2164
$ws.Range("L107").Value = 201.25
$ws.Range("M107").Value = -244
$ws.Range("N107").Value = -4041.25

$ws.Range("H122").Value = 527484.25
$ws.Range("I122").Value = 527484.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1582452.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1580002.75
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2383.3462
$ws.Range("I7").Value = 1832
$ws.Range("J7").Value = 8999.5
$ws.Range("K7").Value = 1832
$ws.Range("L7").Value = 8999.5
$ws.Range("M7").Value = -1720
$ws.Range("N7").Value = -9223.5

$ws.Range("H40").Value = 7827.3213
$ws.Range("I40").Value = 6615.3887
$ws.Range("K40").Value = 6615.3887
$ws.Range("M40").Value = -6479.3887

$ws.Range("H61").Value = 6174282.5
$ws.Range("I61").Value = 7937831.5
$ws.Range("K61").Value = 7937831.5
$ws.Range("M61").Value = -7937629.5

$ws.Range("H113").Value = 6174282.5
$ws.Range("I113").Value = 7937831.5
$ws.Range("K113").Value = 7937831.5
$ws.Range("M113").Value = -7935661.5

$ws.Range("H126").Value = 2383.3462
$ws.Range("I126").Value = 1832
$ws.Range("J126").Value = 8999.5
$ws.Range("K126").Value = 5496
$ws.Range("L126").Value = 26998.5
$ws.Range("M126").Value = -3026
$ws.Range("N126").Value = -31938.5

$ws.Range("H132").Value = 2958.0557
$ws.Range("I132").Value = 2601.2856
$ws.Range("K132").Value = 7803.8568
$ws.Range("M132").Value = -5273.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 889
$ws.Range("I113").Value = 915.5
$ws.Range("J113").Value = 831.5833
$ws.Range("K113").Value = 2746.5
$ws.Range("L113").Value = 2494.7499
$ws.Range("M113").Value = -576.5
$ws.Range("N113").Value = -6834.7499

$ws.Range("H126").Value = 3506.9473
$ws.Range("I126").Value = 3217.2122
$ws.Range("K126").Value = 9651.6366
$ws.Range("M126").Value = -7181.6366

$ws.Range("H132").Value = 18057394
$ws.Range("I132").Value = 22226666
$ws.Range("K132").Value = 66679998
$ws.Range("M132").Value = -66677468

$ws.Range("H136").Value = 5892.3076
$ws.Range("I136").Value = 7145.5
$ws.Range("K136").Value = 21436.5
$ws.Range("M136").Value = -18886.5
